# Card changes for sheet/TrapCardData.xlsx
# "The Trade machine now has all its effect being optional."
#
# The "交换机" (Trade machine) row's effect text is updated so that both of
# its triggered effects become optional ("可选"), and the first effect is
# narrowed from "前方1行或后方1行" (the row in front or the row behind) to
# just "前方" (in front).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newEffect = "回合结束时在房间区，可选：选本牌前方的1张怪物牌，替换房间区或手牌的1张怪物牌。<br>`n回合结束时在手牌区，可选：选手牌的1张怪物牌，替换房间区的1张怪物牌。"

$ws.Range("D15").Value = $newEffect

# Restore the cursor / selection to where the author left it after editing.
$ws.Range("D13").Select()
